$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. Everything that was in B..I shifts to C..J.
$ws.Range("B:B").Insert()

# New column B is a duplicate of the pandas row-index (same values as column A),
# headed "Unnamed: 0".
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Unnamed: 0"

$ws.Range("B2:B22").ClearFormats()
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
$ws.Range("B7").Value = 5
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 7
$ws.Range("B10").Value = 8
$ws.Range("B11").Value = 9
$ws.Range("B12").Value = 10
$ws.Range("B13").Value = 11
$ws.Range("B14").Value = 12
$ws.Range("B15").Value = 13
$ws.Range("B16").Value = 14
$ws.Range("B17").Value = 15
$ws.Range("B18").Value = 16
$ws.Range("B19").Value = 17
$ws.Range("B20").Value = 18
$ws.Range("B21").Value = 19
$ws.Range("B22").Value = 20

# Recomputed heat-map weight values now live in column J (was column I before the
# new column was inserted).
$ws.Range("J2").Value = 0.6605496581398433
$ws.Range("J3").Value = 0.949387713739693
$ws.Range("J4").Value = 0.8316302071898446
$ws.Range("J5").Value = 0.6605496581398433
$ws.Range("J6").Value = 0.949387713739693
$ws.Range("J7").Value = 0.6605496581398433
$ws.Range("J8").Value = 0.8316302071898446
$ws.Range("J9").Value = 0.949387713739693
$ws.Range("J10").Value = 1
$ws.Range("J11").Value = 0.8316302071898446
$ws.Range("J12").Value = 0.8316302071898446
$ws.Range("J13").Value = 0.8316302071898446
$ws.Range("J14").Value = 1
$ws.Range("J15").Value = 0.949387713739693
$ws.Range("J16").Value = 1
$ws.Range("J17").Value = 0.949387713739693
$ws.Range("J18").Value = 1
$ws.Range("J19").Value = 0.6605496581398433
$ws.Range("J20").Value = 0.6605496581398433
$ws.Range("J21").Value = 1
$ws.Range("J22").Value = 0.949387713739693
